# Daily Work Progress update: add a new day-row (row 19) to the log,
# matching the existing formatting used by the previous entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (18) down into the
# new row (19) so the new entry keeps the same date/text/border styling.
$ws.Range("A18:C18").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new day's data.
$ws.Range("A19").Value = 45017   # 01/04/2023
$ws.Range("B19").Value = "Regulator Drafting Program:PLAN Drawing:WingWall"
$ws.Range("C19").Value = "A.K.M Saifuddin"

# Keep selection parked where the user would naturally end up after typing
# the new row (mirrors the workbook's recorded cursor position).
$ws.Range("C23").Select()
